$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.645.98'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.598.43'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.19'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0617'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.46'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0836'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.671.32'
$ws.Range('E12').Value = '  +4.82%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.824.51'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.02'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.521'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.72'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '26.644.75'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '208.03'
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.93'
$ws.Range('E21').Value = '  +2.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.25'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.30'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.84'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.74'
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.14'
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.114'
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.24'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0504'
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.23'
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.654'
$ws.Range('E33').Value = '  -1.42%  '
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('D35').Value = '1.281.89'
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.49'
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.841'
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.46'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.84'
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.919'
$ws.Range('E45').Value = '  +9.33%  '
$ws.Range('D46').Value = '1.736.22'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.66'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.59'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.101'
$ws.Range('E50').Value = '  +3.57%  '
$ws.Range('E51').Value = '  -1.25%  '
